# Add Indian MF 1st Stab
# Inserts 9 new weekly date columns (Jun_16 .. Sep_08) at the front of the
# date-tracking grid (before former column B), shifting the existing weeks
# to the right, and records two new rating-change annotations that fall
# inside the newly added weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 9 new (blank) columns right before the first date column (B),
#    pushing the existing week columns (old B:V) to the right (new K:AE).
#    Excel automatically shifts cell values/formats, so every existing
#    analyst row keeps its data and styling.
$ws.Columns("B:J").Insert()

# 2. Populate the new header row (row 1) with the 9 new week labels, in the
#    same newest-first left-to-right ordering already used by the sheet.
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# 3. The newly inserted cells for every analyst row default to blank; fill
#    them with the same "UN" (unchanged) placeholder used throughout the
#    rest of the grid, matching each row's existing extent.
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $lastCol = $ws.Cells.Item($r, 300).End(-4159).Column()
    if ($lastCol -ge 11) {
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 10)).Value = "UN"
    }
}

# 4. Row 22 (BidaskClub) gets two brand-new rating-change notes that land in
#    the freshly added weeks: an upgrade in the Sep_08 week (col B) and a
#    downgrade in the Jul_07 week (col G). Mark them with the same
#    light-green / light-pink fills used for every other Upgrades /
#    Downgrades note on the sheet.
$upgrade = $ws.Range("B22")
$upgrade.Value = "9/6/2019,Upgrades,Sell -> Hold,"
$upgrade.Interior.Color = 13434828

$downgrade = $ws.Range("G22")
$downgrade.Value = "7/6/2019,Downgrades,Sell -> Strong Sell,"
$downgrade.Interior.Color = 13408767
